# Refresh the cryptocurrency price/volume table with the latest scrape.
# D/E columns hold numeric-looking data stored as TEXT (mirrors the source sheet,
# which persists these as inline strings, not numbers). Excel auto-converts a bare
# "582.33"-style literal to a number, so numeric-looking values are written with a
# leading apostrophe (the standard Excel text-prefix) to keep them as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.125.56'
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("D3").Value = '2.465.42'
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''582.33'
$ws.Range("E5").Value = '  -1.60%  '
$ws.Range("D6").Value = '''167.29'
$ws.Range("E6").Value = '  -3.76%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '''0.515'
$ws.Range("E8").Value = '  -2.03%  '
$ws.Range("D9").Value = '2.466.48'
$ws.Range("E9").Value = '  -1.11%  '
$ws.Range("E10").Value = '  -4.23%  '
$ws.Range("E11").Value = '  -1.11%  '
$ws.Range("D12").Value = '''4.92'
$ws.Range("E12").Value = '  -3.35%  '
$ws.Range("E13").Value = '  -2.70%  '
$ws.Range("D14").Value = '2.931.92'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '''25.46'
$ws.Range("E15").Value = '  -3.13%  '
$ws.Range("D16").Value = '67.031.70'
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("D17").Value = '''0.0000169'
$ws.Range("E17").Value = '  -4.48%  '
$ws.Range("D18").Value = '2.459.67'
$ws.Range("E18").Value = '  -1.62%  '
$ws.Range("D19").Value = '''11.36'
$ws.Range("E19").Value = '  -3.08%  '
$ws.Range("D20").Value = '''7.68'
$ws.Range("E20").Value = '  -4.48%  '
$ws.Range("D21").Value = '''354.59'
$ws.Range("E21").Value = '  -2.84%  '
$ws.Range("D22").Value = '''4.03'
$ws.Range("E22").Value = '  -2.24%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '''69.58'
$ws.Range("E24").Value = '  -2.26%  '
$ws.Range("D25").Value = '''4.23'
$ws.Range("E25").Value = '  -7.24%  '
$ws.Range("E26").Value = '  -7.45%  '
$ws.Range("D27").Value = '''8.97'
$ws.Range("E27").Value = '  -8.99%  '
$ws.Range("D28").Value = '''0.997'
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").Value = '2.589.94'
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("D30").Value = '0.0₃0901'
$ws.Range("E30").Value = '  -6.74%  '
$ws.Range("D31").Value = '''509.99'
$ws.Range("E31").Value = '  -4.24%  '
$ws.Range("D32").Value = '''7.80'
$ws.Range("E32").Value = '  -5.31%  '
$ws.Range("E33").Value = '  -4.43%  '
$ws.Range("D34").Value = '''1.23'
$ws.Range("E34").Value = '  -5.12%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = '''0.119'
$ws.Range("E36").Value = '  -6.82%  '
$ws.Range("D37").Value = '''158.61'
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").Value = '''18.47'
$ws.Range("E38").Value = '  -0.74%  '
$ws.Range("D39").Value = '''18.58'
$ws.Range("E39").Value = '  -0.30%  '
$ws.Range("E40").Value = '  -6.09%  '
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").Value = '''0.328'
$ws.Range("E41").Value = '  -5.90%  '
$ws.Range("D42").Value = '''1.67'
$ws.Range("E42").Value = '  -6.29%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").Value = '''4.79'
$ws.Range("E43").Value = '  -6.24%  '
$ws.Range("B44").Value = 'OKB'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D44").Value = '''38.69'
$ws.Range("E44").Value = '  -2.80%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").Value = '''2.32'
$ws.Range("E45").Value = '  -7.10%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''141.49'
$ws.Range("E46").Value = '  -2.37%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value = '''3.47'
$ws.Range("E47").Value = '  -5.56%  '
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '''0.515'
$ws.Range("E48").Value = '  -6.02%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0254'
$ws.Range("E49").Value = '  -6.76%  '
$ws.Range("B50").Value = 'Optimism'
$ws.Range("C50").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D50").Value = '''1.60'
$ws.Range("E50").Value = '  -5.81%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '''0.0734'
$ws.Range("E51").Value = '  -2.05%  '
